# Update gh-pages output (杭州-漫展信息.xlsx) to the values generated at 456a3b4.
# Column F = "想去人数" (want-to-go count), column G = "最低票价" (lowest price).
# A few events sold out, so their price cell becomes the text "已售罄".

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$f1 = @{3=7593; 7=793; 8=580; 10=66; 13=3044; 14=185; 15=79; 16=712; 17=749; 18=44; 19=443; 20=22; 21=196; 22=210; 23=253; 24=281; 25=124; 27=234; 30=484; 31=440; 32=26; 35=78}
foreach ($r in $f1.Keys) {
    $ws1.Cells.Item($r, 6).Value = $f1[$r]
}
# Rows that sold out: price (G) switches from a number to the text "已售罄"
$ws1.Cells.Item(24, 7).Value = "已售罄"
$ws1.Cells.Item(30, 7).Value = "已售罄"

# ---- Sheet "本地生活" (local life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 202

# ---- Sheet "全部类型" (all types, aggregated view) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$f4 = @{2=202; 6=7593; 10=793; 11=580; 13=66; 17=3044; 18=185; 19=79; 21=712; 22=749; 24=44; 25=443; 26=22; 27=196; 28=210; 29=253; 30=281; 31=124; 33=234; 36=484; 37=440; 38=26; 41=78}
foreach ($r in $f4.Keys) {
    $ws4.Cells.Item($r, 6).Value = $f4[$r]
}
$ws4.Cells.Item(30, 7).Value = "已售罄"
$ws4.Cells.Item(36, 7).Value = "已售罄"
